$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray header row of raw numbers 1..9 at the top (row 1)
$ws.Rows.Item(1).Delete()

# Remove the trailing "\n" column (column I) that held all-zero values
$ws.Columns.Item(9).Delete()

# Update selection to match target (I1:I1048576 selected, as column I was just cleared)
$ws.Range("I1:I1048576").Select() | Out-Null
